$d = $word.ActiveDocument

$d.Content.Find.Execute("18+4=22", $true, $false, $false, $false, $false, $true, 1, $false, "9+14=23", 2)
$d.Content.Find.Execute("54-34=20", $true, $false, $false, $false, $false, $true, 1, $false, "37+53=90", 2)
$d.Content.Find.Execute("99-34=65", $true, $false, $false, $false, $false, $true, 1, $false, "93-27=66", 2)
$d.Content.Find.Execute("89+2=91", $true, $false, $false, $false, $false, $true, 1, $false, "25+38=63", 2)
$d.Content.Find.Execute("32+14=46", $true, $false, $false, $false, $false, $true, 1, $false, "46-7=39", 2)
$d.Content.Find.Execute("23+32=55", $true, $false, $false, $false, $false, $true, 1, $false, "20+17=37", 2)
$d.Content.Find.Execute("97+0=97", $true, $false, $false, $false, $false, $true, 1, $false, "25-8=17", 2)
$d.Content.Find.Execute("49-24=25", $true, $false, $false, $false, $false, $true, 1, $false, "2+67=69", 2)
$d.Content.Find.Execute("11+84=95", $true, $false, $false, $false, $false, $true, 1, $false, "66+10=76", 2)
$d.Content.Find.Execute("32-8=24", $true, $false, $false, $false, $false, $true, 1, $false, "97-65=32", 2)
$d.Content.Find.Execute("46+42=88", $true, $false, $false, $false, $false, $true, 1, $false, "75-14=61", 2)
$d.Content.Find.Execute("19+6=25", $true, $false, $false, $false, $false, $true, 1, $false, "14+44=58", 2)
$d.Content.Find.Execute("76-54=22", $true, $false, $false, $false, $false, $true, 1, $false, "72-56=16", 2)
$d.Content.Find.Execute("32+48=80", $true, $false, $false, $false, $false, $true, 1, $false, "18+33=51", 2)
$d.Content.Find.Execute("38-3=35", $true, $false, $false, $false, $false, $true, 1, $false, "1+29=30", 2)
$d.Content.Find.Execute("87-31=56", $true, $false, $false, $false, $false, $true, 1, $false, "59+27=86", 2)
$d.Content.Find.Execute("70-16=54", $true, $false, $false, $false, $false, $true, 1, $false, "74-35=39", 2)
$d.Content.Find.Execute("22-11=11", $true, $false, $false, $false, $false, $true, 1, $false, "19-10=9", 2)
$d.Content.Find.Execute("61+2=63", $true, $false, $false, $false, $false, $true, 1, $false, "40+54=94", 2)
$d.Content.Find.Execute("94+1=95", $true, $false, $false, $false, $false, $true, 1, $false, "43+8=51", 2)
$d.Content.Find.Execute("3+54=57", $true, $false, $false, $false, $false, $true, 1, $false, "57-41=16", 2)
$d.Content.Find.Execute("25-6=19", $true, $false, $false, $false, $false, $true, 1, $false, "52-26=26", 2)
$d.Content.Find.Execute("17+25=42", $true, $false, $false, $false, $false, $true, 1, $false, "68-51=17", 2)
$d.Content.Find.Execute("72-2=70", $true, $false, $false, $false, $false, $true, 1, $false, "97-89=8", 2)
$d.Content.Find.Execute("17-5=12", $true, $false, $false, $false, $false, $true, 1, $false, "97-12=85", 2)
$d.Content.Find.Execute("54+43=97", $true, $false, $false, $false, $false, $true, 1, $false, "71+4=75", 2)
$d.Content.Find.Execute("80-5=75", $true, $false, $false, $false, $false, $true, 1, $false, "15+81=96", 2)
$d.Content.Find.Execute("2+93=95", $true, $false, $false, $false, $false, $true, 1, $false, "39+23=62", 2)
$d.Content.Find.Execute("18-15=3", $true, $false, $false, $false, $false, $true, 1, $false, "42+36=78", 2)
$d.Content.Find.Execute("71-55=16", $true, $false, $false, $false, $false, $true, 1, $false, "49-42=7", 2)
$d.Content.Find.Execute("25-4=21", $true, $false, $false, $false, $false, $true, 1, $false, "56+4=60", 2)
$d.Content.Find.Execute("74-50=24", $true, $false, $false, $false, $false, $true, 1, $false, "86-67=19", 2)
$d.Content.Find.Execute("96-27=69", $true, $false, $false, $false, $false, $true, 1, $false, "37-32=5", 2)
$d.Content.Find.Execute("73-39=34", $true, $false, $false, $false, $false, $true, 1, $false, "23+39=62", 2)
$d.Content.Find.Execute("52+36=88", $true, $false, $false, $false, $false, $true, 1, $false, "57-10=47", 2)
$d.Content.Find.Execute("84-65=19", $true, $false, $false, $false, $false, $true, 1, $false, "9+87=96", 2)
$d.Content.Find.Execute("92-81=11", $true, $false, $false, $false, $false, $true, 1, $false, "45+37=82", 2)
$d.Content.Find.Execute("48+23=71", $true, $false, $false, $false, $false, $true, 1, $false, "86-16=70", 2)
$d.Content.Find.Execute("1+9=10", $true, $false, $false, $false, $false, $true, 1, $false, "93-21=72", 2)
$d.Content.Find.Execute("17-1=16", $true, $false, $false, $false, $false, $true, 1, $false, "85-65=20", 2)
$d.Content.Find.Execute("27-6=21", $true, $false, $false, $false, $false, $true, 1, $false, "75+19=94", 2)
$d.Content.Find.Execute("64-62=2", $true, $false, $false, $false, $false, $true, 1, $false, "65+24=89", 2)
$d.Content.Find.Execute("93-47=46", $true, $false, $false, $false, $false, $true, 1, $false, "76-5=71", 2)
$d.Content.Find.Execute("77-76=1", $true, $false, $false, $false, $false, $true, 1, $false, "79-53=26", 2)
$d.Content.Find.Execute("31+3=34", $true, $false, $false, $false, $false, $true, 1, $false, "8+16=24", 2)
$d.Content.Find.Execute("4+40=44", $true, $false, $false, $false, $false, $true, 1, $false, "72-64=8", 2)
$d.Content.Find.Execute("70-55=15", $true, $false, $false, $false, $false, $true, 1, $false, "68-50=18", 2)
$d.Content.Find.Execute("25+26=51", $true, $false, $false, $false, $false, $true, 1, $false, "39+46=85", 2)
$d.Content.Find.Execute("6+79=85", $true, $false, $false, $false, $false, $true, 1, $false, "93-52=41", 2)
$d.Content.Find.Execute("91-81=10", $true, $false, $false, $false, $false, $true, 1, $false, "15-14=1", 2)
$d.Content.Find.Execute("64-42=22", $true, $false, $false, $false, $false, $true, 1, $false, "43-40=3", 2)
$d.Content.Find.Execute("82+9=91", $true, $false, $false, $false, $false, $true, 1, $false, "57-46=11", 2)
$d.Content.Find.Execute("52+24=76", $true, $false, $false, $false, $false, $true, 1, $false, "6+41=47", 2)
$d.Content.Find.Execute("48-26=22", $true, $false, $false, $false, $false, $true, 1, $false, "94-71=23", 2)
$d.Content.Find.Execute("69+0=69", $true, $false, $false, $false, $false, $true, 1, $false, "60+37=97", 2)
$d.Content.Find.Execute("41-41=0", $true, $false, $false, $false, $false, $true, 1, $false, "21+71=92", 2)
$d.Content.Find.Execute("63+23=86", $true, $false, $false, $false, $false, $true, 1, $false, "75-2=73", 2)
$d.Content.Find.Execute("90-69=21", $true, $false, $false, $false, $false, $true, 1, $false, "91-76=15", 2)
$d.Content.Find.Execute("92-65=27", $true, $false, $false, $false, $false, $true, 1, $false, "47+19=66", 2)
$d.Content.Find.Execute("69-64=5", $true, $false, $false, $false, $false, $true, 1, $false, "29+56=85", 2)
$d.Content.Find.Execute("92-88=4", $true, $false, $false, $false, $false, $true, 1, $false, "3+91=94", 2)
$d.Content.Find.Execute("50-46=4", $true, $false, $false, $false, $false, $true, 1, $false, "58+24=82", 2)
$d.Content.Find.Execute("99-49=50", $true, $false, $false, $false, $false, $true, 1, $false, "97-64=33", 2)
$d.Content.Find.Execute("22+66=88", $true, $false, $false, $false, $false, $true, 1, $false, "6-3=3", 2)
$d.Content.Find.Execute("17+8=25", $true, $false, $false, $false, $false, $true, 1, $false, "36-33=3", 2)
$d.Content.Find.Execute("27+5=32", $true, $false, $false, $false, $false, $true, 1, $false, "13+1=14", 2)
$d.Content.Find.Execute("12+33=45", $true, $false, $false, $false, $false, $true, 1, $false, "59-13=46", 2)
$d.Content.Find.Execute("57+18=75", $true, $false, $false, $false, $false, $true, 1, $false, "88-48=40", 2)
$d.Content.Find.Execute("13+80=93", $true, $false, $false, $false, $false, $true, 1, $false, "7+41=48", 2)
$d.Content.Find.Execute("61+23=84", $true, $false, $false, $false, $false, $true, 1, $false, "34+31=65", 2)
$d.Content.Find.Execute("23+11=34", $true, $false, $false, $false, $false, $true, 1, $false, "74-26=48", 2)
$d.Content.Find.Execute("25+39=64", $true, $false, $false, $false, $false, $true, 1, $false, "99+0=99", 2)
$d.Content.Find.Execute("98-32=66", $true, $false, $false, $false, $false, $true, 1, $false, "51+7=58", 2)
$d.Content.Find.Execute("2+90=92", $true, $false, $false, $false, $false, $true, 1, $false, "44+31=75", 2)
$d.Content.Find.Execute("4+81=85", $true, $false, $false, $false, $false, $true, 1, $false, "76-0=76", 2)
$d.Content.Find.Execute("16+76=92", $true, $false, $false, $false, $false, $true, 1, $false, "6+33=39", 2)
$d.Content.Find.Execute("74-11=63", $true, $false, $false, $false, $false, $true, 1, $false, "99-22=77", 2)
$d.Content.Find.Execute("98-78=20", $true, $false, $false, $false, $false, $true, 1, $false, "57-51=6", 2)
$d.Content.Find.Execute("2+1=3", $true, $false, $false, $false, $false, $true, 1, $false, "61-39=22", 2)
$d.Content.Find.Execute("94-6=88", $true, $false, $false, $false, $false, $true, 1, $false, "31-25=6", 2)
$d.Content.Find.Execute("85-33=52", $true, $false, $false, $false, $false, $true, 1, $false, "58-34=24", 2)
$d.Content.Find.Execute("27+64=91", $true, $false, $false, $false, $false, $true, 1, $false, "3+21=24", 2)
$d.Content.Find.Execute("76-10=66", $true, $false, $false, $false, $false, $true, 1, $false, "3+26=29", 2)
$d.Content.Find.Execute("25+59=84", $true, $false, $false, $false, $false, $true, 1, $false, "13+61=74", 2)
$d.Content.Find.Execute("30+28=58", $true, $false, $false, $false, $false, $true, 1, $false, "33+25=58", 2)
$d.Content.Find.Execute("65-41=24", $true, $false, $false, $false, $false, $true, 1, $false, "61-48=13", 2)
$d.Content.Find.Execute("47+50=97", $true, $false, $false, $false, $false, $true, 1, $false, "43-11=32", 2)
$d.Content.Find.Execute("57-11=46", $true, $false, $false, $false, $false, $true, 1, $false, "8+81=89", 2)
$d.Content.Find.Execute("82-6=76", $true, $false, $false, $false, $false, $true, 1, $false, "53+11=64", 2)
$d.Content.Find.Execute("47+35=82", $true, $false, $false, $false, $false, $true, 1, $false, "81-5=76", 2)
$d.Content.Find.Execute("79-26=53", $true, $false, $false, $false, $false, $true, 1, $false, "97-52=45", 2)
$d.Content.Find.Execute("60-13=47", $true, $false, $false, $false, $false, $true, 1, $false, "42-41=1", 2)
$d.Content.Find.Execute("68-44=24", $true, $false, $false, $false, $false, $true, 1, $false, "32+52=84", 2)
$d.Content.Find.Execute("24+35=59", $true, $false, $false, $false, $false, $true, 1, $false, "78-2=76", 2)
$d.Content.Find.Execute("59+32=91", $true, $false, $false, $false, $false, $true, 1, $false, "38-18=20", 2)
$d.Content.Find.Execute("8-2=6", $true, $false, $false, $false, $false, $true, 1, $false, "7+8=15", 2)
$d.Content.Find.Execute("33-22=11", $true, $false, $false, $false, $false, $true, 1, $false, "89-46=43", 2)
$d.Content.Find.Execute("61+9=70", $true, $false, $false, $false, $false, $true, 1, $false, "19+13=32", 2)
$d.Content.Find.Execute("99-55=44", $true, $false, $false, $false, $false, $true, 1, $false, "35-13=22", 2)
$d.Content.Find.Execute("54-38=16", $true, $false, $false, $false, $false, $true, 1, $false, "50+35=85", 2)
